# Generate Report for Handback
# Reorders the localization status rows (60d3e31d.. and fd123f5e.. move up,
# marked as "Handed back: in sync with en-US") and fills in the newly
# populated "Latest Target File" / "Latest Handback File" columns on the
# per-locale sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# URL lookup table for hyperlinks, keyed by the display text (the file
# name shown in the cell). These are the same targets already present
# in the workbook; we simply re-attach them (and reuse them for the new
# "Latest Target File"/"Latest Handback File" cells that point at the
# same files as "Source File Name"/"Latest Handoff File").
# ---------------------------------------------------------------------
$urls = @{
    "70058cc9-c335-4816-b6d1-eb41c2e484ea.md" = "https://github.com/OpenLocalizationTest/oltest/blob/0c081f304363a8d5877c53e5daff058cb20a9522/e2e/70058cc9-c335-4816-b6d1-eb41c2e484ea.md";
    "4dd38fef-441a-4be0-b79a-ef6c9247ebcf.md" = "https://github.com/OpenLocalizationTest/oltest/blob/a22127a5bf27b2a36c8d9832d09dc3ecb34c3040/e2e/4dd38fef-441a-4be0-b79a-ef6c9247ebcf.md";
    "60d3e31d-5338-4519-9c0f-1441a4cf07a9.md" = "https://github.com/OpenLocalizationTest/oltest/blob/a22127a5bf27b2a36c8d9832d09dc3ecb34c3040/e2e/60d3e31d-5338-4519-9c0f-1441a4cf07a9.md";
    "fd123f5e-76f2-4b36-91d0-7e2a6043f1a2.md" = "https://github.com/OpenLocalizationTest/oltest/blob/a22127a5bf27b2a36c8d9832d09dc3ecb34c3040/e2e/fd123f5e-76f2-4b36-91d0-7e2a6043f1a2.md";
    ".localization-config" = "https://github.com/OpenLocalizationTest/oltest/blob/a22127a5bf27b2a36c8d9832d09dc3ecb34c3040/.localization-config";

    "70058cc9-c335-4816-b6d1-eb41c2e484ea.d5899fcb1515f857962642b1ceab8a68295ec2fb.zh-cn.xlf" = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b3d92526e38c3918ca1ab380f3756dad2d6c95cb/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/70058cc9-c335-4816-b6d1-eb41c2e484ea.d5899fcb1515f857962642b1ceab8a68295ec2fb.zh-cn.xlf";
    "4dd38fef-441a-4be0-b79a-ef6c9247ebcf.9ebd15eee1cc650407d011344150e433768ce247.zh-cn.xlf" = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0c16320d47a20331fd72a6dca1dfdb8e2fb383d4/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/4dd38fef-441a-4be0-b79a-ef6c9247ebcf.9ebd15eee1cc650407d011344150e433768ce247.zh-cn.xlf";
    "60d3e31d-5338-4519-9c0f-1441a4cf07a9.b0c4d9c873c6d5f3344a83e299fa2f94467cb296.zh-cn.xlf" = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0c16320d47a20331fd72a6dca1dfdb8e2fb383d4/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/60d3e31d-5338-4519-9c0f-1441a4cf07a9.b0c4d9c873c6d5f3344a83e299fa2f94467cb296.zh-cn.xlf";
    "fd123f5e-76f2-4b36-91d0-7e2a6043f1a2.8c47fac789840fd5322964296750e25cfd30b937.zh-cn.xlf" = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0c16320d47a20331fd72a6dca1dfdb8e2fb383d4/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/fd123f5e-76f2-4b36-91d0-7e2a6043f1a2.8c47fac789840fd5322964296750e25cfd30b937.zh-cn.xlf";

    "70058cc9-c335-4816-b6d1-eb41c2e484ea.d5899fcb1515f857962642b1ceab8a68295ec2fb.de-de.xlf" = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8742ffd895306e42a64dc16fd15d18efc2bc2379/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/70058cc9-c335-4816-b6d1-eb41c2e484ea.d5899fcb1515f857962642b1ceab8a68295ec2fb.de-de.xlf";
    "4dd38fef-441a-4be0-b79a-ef6c9247ebcf.9ebd15eee1cc650407d011344150e433768ce247.de-de.xlf" = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/66a3dc743f1b47375809f62342b4bd6e817fdb0e/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/4dd38fef-441a-4be0-b79a-ef6c9247ebcf.9ebd15eee1cc650407d011344150e433768ce247.de-de.xlf";
    "60d3e31d-5338-4519-9c0f-1441a4cf07a9.b0c4d9c873c6d5f3344a83e299fa2f94467cb296.de-de.xlf" = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/66a3dc743f1b47375809f62342b4bd6e817fdb0e/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/60d3e31d-5338-4519-9c0f-1441a4cf07a9.b0c4d9c873c6d5f3344a83e299fa2f94467cb296.de-de.xlf";
    "fd123f5e-76f2-4b36-91d0-7e2a6043f1a2.8c47fac789840fd5322964296750e25cfd30b937.de-de.xlf" = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/66a3dc743f1b47375809f62342b4bd6e817fdb0e/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/fd123f5e-76f2-4b36-91d0-7e2a6043f1a2.8c47fac789840fd5322964296750e25cfd30b937.de-de.xlf";
}

function Set-LinkedCell($ws, $cellRef, $displayText) {
    $ws.Range($cellRef).Value2 = $displayText
    $ws.Hyperlinks.Add($ws.Range($cellRef), $urls[$displayText], "", "", $displayText) | Out-Null
}

# =======================================================================
# Sheet "Overview"
# =======================================================================
$ov = $wb.Worksheets.Item("Overview")
$ov.Hyperlinks.Delete()

$ov.Range("A2").Value2 = "60d3e31d-5338-4519-9c0f-1441a4cf07a9.md"
$ov.Range("B2").Value2 = "Handed back: in sync with en-US"
$ov.Range("C2").Value2 = "Handed back: in sync with en-US"

$ov.Range("A3").Value2 = "fd123f5e-76f2-4b36-91d0-7e2a6043f1a2.md"
$ov.Range("B3").Value2 = "Handed back: in sync with en-US"
$ov.Range("C3").Value2 = "Handed back: in sync with en-US"

$ov.Range("A4").Value2 = "70058cc9-c335-4816-b6d1-eb41c2e484ea.md"
$ov.Range("B4").Value2 = "In Translation"
$ov.Range("C4").Value2 = "In Translation"

$ov.Range("A5").Value2 = "4dd38fef-441a-4be0-b79a-ef6c9247ebcf.md"
$ov.Range("B5").Value2 = "Ready for handoff"
$ov.Range("C5").Value2 = "Ready for handoff"

$ov.Range("A6").Value2 = ".localization-config"
$ov.Range("B6").Value2 = "Not to be localized"
$ov.Range("C6").Value2 = "Not to be localized"

Set-LinkedCell $ov "A2" "60d3e31d-5338-4519-9c0f-1441a4cf07a9.md"
Set-LinkedCell $ov "A3" "fd123f5e-76f2-4b36-91d0-7e2a6043f1a2.md"
Set-LinkedCell $ov "A4" "70058cc9-c335-4816-b6d1-eb41c2e484ea.md"
Set-LinkedCell $ov "A5" "4dd38fef-441a-4be0-b79a-ef6c9247ebcf.md"
Set-LinkedCell $ov "A6" ".localization-config"

# =======================================================================
# Sheet "zh-cn"
# =======================================================================
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Hyperlinks.Delete()

# Row 2 : 60d3e31d...
$zh.Range("A2").Value2 = "60d3e31d-5338-4519-9c0f-1441a4cf07a9.md"
$zh.Range("B2").Value2 = "Handed back: in sync with en-US"
$zh.Range("C2").Value2 = "60d3e31d-5338-4519-9c0f-1441a4cf07a9.b0c4d9c873c6d5f3344a83e299fa2f94467cb296.zh-cn.xlf"
$zh.Range("D2").Value2 = "2016-03-10 04:15:40"
$zh.Range("E2").Value2 = "60d3e31d-5338-4519-9c0f-1441a4cf07a9.md"
$zh.Range("F2").Value2 = "60d3e31d-5338-4519-9c0f-1441a4cf07a9.b0c4d9c873c6d5f3344a83e299fa2f94467cb296.zh-cn.xlf"
$zh.Range("G2").Value2 = "2016-03-10 04:16:33"
$zh.Range("H2").Value2 = "Include"

# Row 3 : fd123f5e...
$zh.Range("A3").Value2 = "fd123f5e-76f2-4b36-91d0-7e2a6043f1a2.md"
$zh.Range("B3").Value2 = "Handed back: in sync with en-US"
$zh.Range("C3").Value2 = "fd123f5e-76f2-4b36-91d0-7e2a6043f1a2.8c47fac789840fd5322964296750e25cfd30b937.zh-cn.xlf"
$zh.Range("D3").Value2 = "2016-03-10 04:15:40"
$zh.Range("E3").Value2 = "fd123f5e-76f2-4b36-91d0-7e2a6043f1a2.md"
$zh.Range("F3").Value2 = "fd123f5e-76f2-4b36-91d0-7e2a6043f1a2.8c47fac789840fd5322964296750e25cfd30b937.zh-cn.xlf"
$zh.Range("G3").Value2 = "2016-03-10 04:16:33"
$zh.Range("H3").Value2 = "Include"

# Row 4 : 70058cc9...
$zh.Range("A4").Value2 = "70058cc9-c335-4816-b6d1-eb41c2e484ea.md"
$zh.Range("B4").Value2 = "In Translation"
$zh.Range("C4").Value2 = "70058cc9-c335-4816-b6d1-eb41c2e484ea.d5899fcb1515f857962642b1ceab8a68295ec2fb.zh-cn.xlf"
$zh.Range("D4").Value2 = "2016-03-10 04:13:01"
$zh.Range("E4").ClearContents()
$zh.Range("F4").ClearContents()
$zh.Range("G4").Value2 = "0001-01-01 00:00:00"
$zh.Range("H4").Value2 = "Include"

# Row 5 : 4dd38fef...
$zh.Range("A5").Value2 = "4dd38fef-441a-4be0-b79a-ef6c9247ebcf.md"
$zh.Range("B5").Value2 = "Ready for handoff"
$zh.Range("C5").Value2 = "4dd38fef-441a-4be0-b79a-ef6c9247ebcf.9ebd15eee1cc650407d011344150e433768ce247.zh-cn.xlf"
$zh.Range("D5").Value2 = "2016-03-10 04:15:40"
$zh.Range("E5").ClearContents()
$zh.Range("F5").ClearContents()
$zh.Range("G5").Value2 = "0001-01-01 00:00:00"
$zh.Range("H5").Value2 = "Include"

# Row 6 : .localization-config
$zh.Range("A6").Value2 = ".localization-config"
$zh.Range("B6").Value2 = "Not to be localized"
$zh.Range("D6").Value2 = "0001-01-01 00:00:00"
$zh.Range("G6").Value2 = "0001-01-01 00:00:00"
$zh.Range("H6").Value2 = "Ignored"

Set-LinkedCell $zh "A2" "60d3e31d-5338-4519-9c0f-1441a4cf07a9.md"
Set-LinkedCell $zh "C2" "60d3e31d-5338-4519-9c0f-1441a4cf07a9.b0c4d9c873c6d5f3344a83e299fa2f94467cb296.zh-cn.xlf"
Set-LinkedCell $zh "E2" "60d3e31d-5338-4519-9c0f-1441a4cf07a9.md"
Set-LinkedCell $zh "F2" "60d3e31d-5338-4519-9c0f-1441a4cf07a9.b0c4d9c873c6d5f3344a83e299fa2f94467cb296.zh-cn.xlf"

Set-LinkedCell $zh "A3" "fd123f5e-76f2-4b36-91d0-7e2a6043f1a2.md"
Set-LinkedCell $zh "C3" "fd123f5e-76f2-4b36-91d0-7e2a6043f1a2.8c47fac789840fd5322964296750e25cfd30b937.zh-cn.xlf"
Set-LinkedCell $zh "E3" "fd123f5e-76f2-4b36-91d0-7e2a6043f1a2.md"
Set-LinkedCell $zh "F3" "fd123f5e-76f2-4b36-91d0-7e2a6043f1a2.8c47fac789840fd5322964296750e25cfd30b937.zh-cn.xlf"

Set-LinkedCell $zh "A4" "70058cc9-c335-4816-b6d1-eb41c2e484ea.md"
Set-LinkedCell $zh "C4" "70058cc9-c335-4816-b6d1-eb41c2e484ea.d5899fcb1515f857962642b1ceab8a68295ec2fb.zh-cn.xlf"

Set-LinkedCell $zh "A5" "4dd38fef-441a-4be0-b79a-ef6c9247ebcf.md"
Set-LinkedCell $zh "C5" "4dd38fef-441a-4be0-b79a-ef6c9247ebcf.9ebd15eee1cc650407d011344150e433768ce247.zh-cn.xlf"

Set-LinkedCell $zh "A6" ".localization-config"

# =======================================================================
# Sheet "de-de"
# =======================================================================
$de = $wb.Worksheets.Item("de-de")
$de.Hyperlinks.Delete()

# Row 2 : 60d3e31d...
$de.Range("A2").Value2 = "60d3e31d-5338-4519-9c0f-1441a4cf07a9.md"
$de.Range("B2").Value2 = "Handed back: in sync with en-US"
$de.Range("C2").Value2 = "60d3e31d-5338-4519-9c0f-1441a4cf07a9.b0c4d9c873c6d5f3344a83e299fa2f94467cb296.de-de.xlf"
$de.Range("D2").Value2 = "2016-03-10 04:15:44"
$de.Range("E2").Value2 = "60d3e31d-5338-4519-9c0f-1441a4cf07a9.md"
$de.Range("F2").Value2 = "60d3e31d-5338-4519-9c0f-1441a4cf07a9.b0c4d9c873c6d5f3344a83e299fa2f94467cb296.de-de.xlf"
$de.Range("G2").Value2 = "2016-03-10 04:16:40"
$de.Range("H2").Value2 = "Include"

# Row 3 : fd123f5e...
$de.Range("A3").Value2 = "fd123f5e-76f2-4b36-91d0-7e2a6043f1a2.md"
$de.Range("B3").Value2 = "Handed back: in sync with en-US"
$de.Range("C3").Value2 = "fd123f5e-76f2-4b36-91d0-7e2a6043f1a2.8c47fac789840fd5322964296750e25cfd30b937.de-de.xlf"
$de.Range("D3").Value2 = "2016-03-10 04:15:44"
$de.Range("E3").Value2 = "fd123f5e-76f2-4b36-91d0-7e2a6043f1a2.md"
$de.Range("F3").Value2 = "fd123f5e-76f2-4b36-91d0-7e2a6043f1a2.8c47fac789840fd5322964296750e25cfd30b937.de-de.xlf"
$de.Range("G3").Value2 = "2016-03-10 04:16:40"
$de.Range("H3").Value2 = "Include"

# Row 4 : 70058cc9...
$de.Range("A4").Value2 = "70058cc9-c335-4816-b6d1-eb41c2e484ea.md"
$de.Range("B4").Value2 = "In Translation"
$de.Range("C4").Value2 = "70058cc9-c335-4816-b6d1-eb41c2e484ea.d5899fcb1515f857962642b1ceab8a68295ec2fb.de-de.xlf"
$de.Range("D4").Value2 = "2016-03-10 04:13:23"
$de.Range("E4").ClearContents()
$de.Range("F4").ClearContents()
$de.Range("G4").Value2 = "0001-01-01 00:00:00"
$de.Range("H4").Value2 = "Include"

# Row 5 : 4dd38fef...
$de.Range("A5").Value2 = "4dd38fef-441a-4be0-b79a-ef6c9247ebcf.md"
$de.Range("B5").Value2 = "Ready for handoff"
$de.Range("C5").Value2 = "4dd38fef-441a-4be0-b79a-ef6c9247ebcf.9ebd15eee1cc650407d011344150e433768ce247.de-de.xlf"
$de.Range("D5").Value2 = "2016-03-10 04:15:44"
$de.Range("E5").ClearContents()
$de.Range("F5").ClearContents()
$de.Range("G5").Value2 = "0001-01-01 00:00:00"
$de.Range("H5").Value2 = "Include"

# Row 6 : .localization-config
$de.Range("A6").Value2 = ".localization-config"
$de.Range("B6").Value2 = "Not to be localized"
$de.Range("D6").Value2 = "0001-01-01 00:00:00"
$de.Range("G6").Value2 = "0001-01-01 00:00:00"
$de.Range("H6").Value2 = "Ignored"

Set-LinkedCell $de "A2" "60d3e31d-5338-4519-9c0f-1441a4cf07a9.md"
Set-LinkedCell $de "C2" "60d3e31d-5338-4519-9c0f-1441a4cf07a9.b0c4d9c873c6d5f3344a83e299fa2f94467cb296.de-de.xlf"
Set-LinkedCell $de "E2" "60d3e31d-5338-4519-9c0f-1441a4cf07a9.md"
Set-LinkedCell $de "F2" "60d3e31d-5338-4519-9c0f-1441a4cf07a9.b0c4d9c873c6d5f3344a83e299fa2f94467cb296.de-de.xlf"

Set-LinkedCell $de "A3" "fd123f5e-76f2-4b36-91d0-7e2a6043f1a2.md"
Set-LinkedCell $de "C3" "fd123f5e-76f2-4b36-91d0-7e2a6043f1a2.8c47fac789840fd5322964296750e25cfd30b937.de-de.xlf"
Set-LinkedCell $de "E3" "fd123f5e-76f2-4b36-91d0-7e2a6043f1a2.md"
Set-LinkedCell $de "F3" "fd123f5e-76f2-4b36-91d0-7e2a6043f1a2.8c47fac789840fd5322964296750e25cfd30b937.de-de.xlf"

Set-LinkedCell $de "A4" "70058cc9-c335-4816-b6d1-eb41c2e484ea.md"
Set-LinkedCell $de "C4" "70058cc9-c335-4816-b6d1-eb41c2e484ea.d5899fcb1515f857962642b1ceab8a68295ec2fb.de-de.xlf"

Set-LinkedCell $de "A5" "4dd38fef-441a-4be0-b79a-ef6c9247ebcf.md"
Set-LinkedCell $de "C5" "4dd38fef-441a-4be0-b79a-ef6c9247ebcf.9ebd15eee1cc650407d011344150e433768ce247.de-de.xlf"

Set-LinkedCell $de "A6" ".localization-config"

Write-Host "Report regenerated for handback."
